$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Sheet1")

# Route B stops (row 2) renamed to reflect dispatch-authority direction suffixes
$ws.Range("B2").Value = "Dormont-N"
$ws.Range("E2").Value = "Glenbury-U"
$ws.Range("F2").Value = "Overbrook-W"
$ws.Range("G2").Value = "Central-W"

# Route D stops (row 4)
$ws.Range("B4").Value = "Central-I"
$ws.Range("C4").Value = "Inglewood-I"

# Route E stops (row 5)
$ws.Range("D5").Value = "Overbrook-I"
$ws.Range("E5").Value = "Central-W"

$ws.Range("G6").Select()
